$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Talent_Acquisition")
$ws1.Range("BC1:BF1").Style = "Normal"
$ws1.Range("A8:BC8").Style = "Normal"
